$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J column (J2:J11)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New bold style cells with labels + stats in rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Apply bold font, size 12, vertical-center alignment to the B14 cell first,
# then propagate the format to B15:B17 via a format-only copy/paste so the
# style table doesn't accumulate unused intermediate entries.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108  # xlVAlignCenter

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A14:A17").RowHeight = 15.6

# Page setup (Letter-ish A4/paper size 9 = A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A14:B17").Select()
